$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$sh = $s.Shapes.Item(2)
$tf = $sh.TextFrame
$tr = $tf.TextRange

# Turn on "Shrink text on overflow" (normAutofit) for the body text.
$tf.AutoSize = 2

# Insert a brand-new first paragraph ("Implemented in C++ and Qt") before
# the existing "Kactus2 goals" paragraph.
$firstPara = $tr.Paragraphs(1)
$firstPara.InsertBefore("Implemented in C++ and Qt" + [char]13)

# Re-fetch that freshly created paragraph and mark it as Finnish, matching
# the language used elsewhere for short, ad-hoc phrases in this deck.
$impPara = $tr.Paragraphs(1)
$impPara.LanguageID = "fi-FI"

# "Kactus2 goals" -> "Kactus2 " + "goals"
$kactusPara = $tr.Paragraphs(2)
$kactusFirst = $kactusPara.Characters(1, 8)
Write-Host $kactusFirst.Text

# "Handling file dependencies is essential for good usability" ->
#   "Handling " + "file dependencies is essential for good usability"
$handlingPara = $tr.Paragraphs(6)
$handlingFirst = $handlingPara.Characters(1, 9)
Write-Host $handlingFirst.Text

Write-Host $tr.Text
